$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Updated country rankings/case data (re-sorted by Casos totales desc).
# Each entry: spreadsheet row number, country name (col A), and the
# 7 numeric stat columns B:H -> Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes.
$updates = @(
    @{ Row = 4; Name = "Estados Unidos"; Values = @(686431, 8861, 58179, 592674, 13380, 961, 35578) }
    @{ Row = 8; Name = "Alemania"; Values = @(138497, 799, 81800, 52504, 4288, 141, 4193) }
    @{ Row = 20; Name = "Austria"; Values = @(14579, 103, 9704, 4465, 227, 0, 410) }
    @{ Row = 24; Name = "Israel"; Values = @(12982, 224, 3126, 9705, 168, 9, 151) }
    @{ Row = 57; Name = "Marruecos"; Values = @(2564, 281, 281, 2148, 1, 5, 135) }
    @{ Row = 100; Name = "Guinea"; Values = @(477, 39, 59, 415, 0, 2, 3) }
    @{ Row = 101; Name = "Bolivia"; Values = @(465, 24, 26, 408, 3, 2, 31) }
    @{ Row = 102; Name = "Honduras"; Values = @(442, 16, 10, 391, 10, 6, 41) }
    @{ Row = 103; Name = "Nigeria"; Values = @(442, 0, 152, 277, 2, 0, 13) }
    @{ Row = 106; Name = "Jordania"; Values = @(407, 5, 259, 141, 5, 0, 7) }
    @{ Row = 107; Name = "Estado de Palestina"; Values = @(402, 28, 69, 331, 0, 0, 2) }
    @{ Row = 108; Name = "Reunion"; Values = @(402, 8, 237, 165, 4, 0, 0) }
    @{ Row = 109; Name = "Taiwan"; Values = @(395, 0, 166, 223, 0, 0, 6) }
    @{ Row = 119; Name = "Sri Lanka"; Values = @(244, 6, 77, 160, 1, 0, 7) }
    @{ Row = 163; Name = "Siria"; Values = @(38, 5, 5, 31, 0, 0, 2) }
    @{ Row = 164; Name = "Eritrea"; Values = @(35, 0, 0, 35, 0, 0, 0) }
    @{ Row = 165; Name = "San Martin (Parte Francesa)"; Values = @(35, 0, 13, 20, 5, 0, 2) }
    @{ Row = 166; Name = "Benin"; Values = @(35, 0, 18, 16, 0, 0, 1) }
    @{ Row = 167; Name = "Mozambique"; Values = @(34, 3, 2, 32, 0, 0, 0) }
    @{ Row = 190; Name = "San Cristobal y Nieves"; Values = @(14, 0, 0, 14, 0, 0, 0) }
    @{ Row = 191; Name = "Granada"; Values = @(14, 0, 0, 14, 2, 0, 0) }
    @{ Row = 195; Name = "Montserrat"; Values = @(11, 0, 1, 10, 1, 0, 0) }
    @{ Row = 196; Name = "Islas Malvinas"; Values = @(11, 0, 1, 10, 0, 0, 0) }
    @{ Row = 197; Name = "Islas Turcas y Caicos"; Values = @(11, 0, 0, 10, 0, 0, 1) }
    @{ Row = 215; Name = "Yemen"; Values = @(1, 0, 0, 1, 0, 0, 0) }
    @{ Row = 216; Name = "San Pedro y Miquelon"; Values = @(1, 0, 0, 1, 0, 0, 0) }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 1).Value = $u.Name
    $vals = $u.Values
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, 2 + $i).Value = $vals[$i]
    }
}
